$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 29: Nro=28, Actor="TarjetaABM" (same as row28), Nombre="traerTarjetaConBeneficios"
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = "TarjetaABM"
$ws.Range("C29").Value = "traerTarjetaConBeneficios"

# Move selection to the newly added row, mirroring the saved workbook state
$ws.Range("D29").Select()
